$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '26.229.38'
$ws.Cells.Item(2, 5).Value = '  -0.55%  '
$ws.Cells.Item(3, 4).Value = '1.588.70'
$ws.Cells.Item(3, 5).Value = '  -0.20%  '
$ws.Cells.Item(4, 4).Value = "'1.00"
$ws.Cells.Item(4, 4).Style = "Normal"
$ws.Cells.Item(4, 5).Value = '  -0.09%  '
$ws.Cells.Item(5, 4).Value = "'211.94"
$ws.Cells.Item(5, 4).Style = "Normal"
$ws.Cells.Item(5, 5).Value = '  +0.81%  '
$ws.Cells.Item(6, 5).Value = '  -0.27%  '
$ws.Cells.Item(7, 5).Value = '  -0.06%  '
$ws.Cells.Item(8, 5).Value = '  -0.33%  '
$ws.Cells.Item(9, 5).Value = '  -1.08%  '
$ws.Cells.Item(10, 5).Value = '  -1.68%  '
$ws.Cells.Item(11, 5).Value = '  +0.23%  '
$ws.Cells.Item(12, 4).Value = '1.812.11'
$ws.Cells.Item(12, 5).Value = '  -0.28%  '
$ws.Cells.Item(13, 4).Value = '1.589.57'
$ws.Cells.Item(13, 5).Value = '  +1.16%  '
$ws.Cells.Item(14, 5).Value = '  -1.50%  '
$ws.Cells.Item(15, 5).Value = '  -0.93%  '
$ws.Cells.Item(16, 5).Value = '  -0.89%  '
$ws.Cells.Item(17, 4).Value = '26.235.78'
$ws.Cells.Item(18, 5).Value = '  -0.33%  '
$ws.Cells.Item(19, 4).Value = "'7.47"
$ws.Cells.Item(19, 4).Style = "Normal"
$ws.Cells.Item(19, 5).Value = '  -0.12%  '
$ws.Cells.Item(20, 4).Value = "'214.13"
$ws.Cells.Item(20, 4).Style = "Normal"
$ws.Cells.Item(20, 5).Value = '  +1.46%  '
$ws.Cells.Item(21, 4).Value = "'1.00"
$ws.Cells.Item(21, 4).Style = "Normal"
$ws.Cells.Item(21, 5).Value = '  -0.03%  '
$ws.Cells.Item(22, 5).Value = '  -0.63%  '
$ws.Cells.Item(23, 4).Value = "'8.97"
$ws.Cells.Item(23, 4).Style = "Normal"
$ws.Cells.Item(23, 5).Value = '  +0.41%  '
$ws.Cells.Item(24, 5).Value = '  -1.48%  '
$ws.Cells.Item(25, 4).Value = "'144.41"
$ws.Cells.Item(25, 4).Style = "Normal"
$ws.Cells.Item(25, 5).Value = '  -0.40%  '
$ws.Cells.Item(26, 5).Value = '  -0.08%  '
$ws.Cells.Item(27, 4).Value = "'6.98"
$ws.Cells.Item(27, 4).Style = "Normal"
$ws.Cells.Item(28, 5).Value = '  -0.94%  '
$ws.Cells.Item(29, 5).Value = '  -1.21%  '
$ws.Cells.Item(30, 5).Value = '  -2.06%  '
$ws.Cells.Item(31, 5).Value = '  +0.25%  '
$ws.Cells.Item(32, 5).Value = '  -1.11%  '
$ws.Cells.Item(33, 4).Value = '1.416.62'
$ws.Cells.Item(33, 5).Value = '  +8.12%  '
$ws.Cells.Item(35, 5).Value = '  -0.54%  '
$ws.Cells.Item(36, 4).Value = "'0.589"
$ws.Cells.Item(36, 4).Style = "Normal"
$ws.Cells.Item(36, 5).Value = '  -4.34%  '
$ws.Cells.Item(37, 4).Value = "'1.45"
$ws.Cells.Item(37, 4).Style = "Normal"
$ws.Cells.Item(37, 5).Value = '  -1.53%  '
$ws.Cells.Item(38, 5).Value = '  -1.32%  '
$ws.Cells.Item(39, 4).Value = "'5.90"
$ws.Cells.Item(39, 4).Style = "Normal"
$ws.Cells.Item(39, 5).Value = '  +4.90%  '
$ws.Cells.Item(40, 5).Value = '  +1.30%  '
$ws.Cells.Item(41, 5).Value = '  -0.07%  '
$ws.Cells.Item(42, 4).Value = "'0.938"
$ws.Cells.Item(42, 4).Style = "Normal"
$ws.Cells.Item(42, 5).Value = '  -14.80%  '
$ws.Cells.Item(43, 4).Value = "'0.765"
$ws.Cells.Item(43, 4).Style = "Normal"
$ws.Cells.Item(43, 5).Value = '  -0.22%  '
$ws.Cells.Item(45, 4).Value = '1.723.52'
$ws.Cells.Item(45, 5).Value = '  -0.31%  '
$ws.Cells.Item(46, 5).Value = '  -2.47%  '
$ws.Cells.Item(47, 4).Value = "'85.84"
$ws.Cells.Item(47, 4).Style = "Normal"
$ws.Cells.Item(47, 5).Value = '  -2.35%  '
$ws.Cells.Item(48, 5).Value = '  -0.67%  '
$ws.Cells.Item(49, 5).Value = '  -0.42%  '
$ws.Cells.Item(50, 4).Value = "'0.0968"
$ws.Cells.Item(50, 4).Style = "Normal"
$ws.Cells.Item(50, 5).Value = '  -1.44%  '
$ws.Cells.Item(51, 4).Value = "'0.998"
$ws.Cells.Item(51, 4).Style = "Normal"
$ws.Cells.Item(51, 5).Value = '  -0.06%  '
